$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.709.05'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '3.079.36'
$ws.Range('E3').Value = '  -2.35%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.93'
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.91'
$ws.Range('E6').Value = '  +3.84%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.533'
$ws.Range('E8').Value = '  +0.27%  '
$ws.Range('D9').Value = '3.076.25'
$ws.Range('E9').Value = '  -2.39%  '
$ws.Range('E10').Value = '  -2.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.92'
$ws.Range('E11').Value = '  -0.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.458'
$ws.Range('E12').Value = '  -0.87%  '
$ws.Range('E13').Value = '  -2.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.66'
$ws.Range('E14').Value = '  +1.23%  '
$ws.Range('D15').Value = '3.591.25'
$ws.Range('E15').Value = '  -2.32%  '
$ws.Range('E16').Value = '  -2.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.17'
$ws.Range('E17').Value = '  -0.84%  '
$ws.Range('D18').Value = '63.604.65'
$ws.Range('E18').Value = '  -0.85%  '
$ws.Range('D19').Value = '3.077.66'
$ws.Range('E19').Value = '  -2.44%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '469.58'
$ws.Range('E20').Value = '  +0.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.64'
$ws.Range('E21').Value = '  +1.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.727'
$ws.Range('E22').Value = '  -1.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.55'
$ws.Range('E23').Value = '  +0.67%  '
$ws.Range('B24').Value = 'Fetch.AI'
$ws.Range('C24').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.38'
$ws.Range('E24').Value = '  +1.65%  '
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.25'
$ws.Range('E25').Value = '  +1.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '81.21'
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.87'
$ws.Range('E28').Value = '  +1.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.69'
$ws.Range('E29').Value = '  -1.35%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.31'
$ws.Range('E30').Value = '  +0.66%  '
$ws.Range('B31').Value = 'FirstDigitalUSD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('E32').Value = '  -1.09%  '
$ws.Range('E33').Value = '  +4.50%  '
$ws.Range('E34').Value = '  -0.71%  '
$ws.Range('D35').Value = '0.0₃0851'
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('E36').Value = '  -1.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.41'
$ws.Range('E37').Value = '  +3.83%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.12'
$ws.Range('E38').Value = '  -0.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.25'
$ws.Range('E39').Value = '  -3.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.34'
$ws.Range('E40').Value = '  +2.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.70'
$ws.Range('E41').Value = '  -2.54%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '448.69'
$ws.Range('E42').Value = '  -1.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.288'
$ws.Range('E43').Value = '  -1.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0366'
$ws.Range('E44').Value = '  -1.83%  '
$ws.Range('B45').Value = 'Arweave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '40.23'
$ws.Range('E45').Value = '  -0.38%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '2.830.52'
$ws.Range('E46').Value = '  -3.37%  '
$ws.Range('E47').Value = '  -0.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.20'
$ws.Range('E48').Value = '  +1.11%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.50'
$ws.Range('E49').Value = '  +4.29%  '
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.26'
$ws.Range('E51').Value = '  +0.58%  '
